# Rerunning southern bluefin tuna models
# Adds a new "SBFT" taxon row to both the "optimal models" and
# "relative importance" sheets.

$wb = $excel.ActiveWorkbook

$wsModels = $wb.Worksheets.Item("optimal models")
$wsImportance = $wb.Worksheets.Item("relative importance")

# --- "relative importance" sheet: append row 12 -------------------------
$wsImportance.Range("A12").Value = "SBFT"
$wsImportance.Range("B12").Value = 9.7753189
$wsImportance.Range("C12").Value = 1.3711376
$wsImportance.Range("D12").Value = 2.0247563
$wsImportance.Range("E12").Value = 0.9137509
$wsImportance.Range("F12").Value = 0.9702147
$wsImportance.Range("G12").Value = 5.4405873
$wsImportance.Range("H12").Value = 7.380599
$wsImportance.Range("I12").Value = 1.9768465
$wsImportance.Range("J12").Value = 1.4763375
$wsImportance.Range("K12").Value = 0.16223
$wsImportance.Range("L12").Value = 2.5284049
$wsImportance.Range("M12").Value = 4.8537287
$wsImportance.Range("N12").Value = 42.1807406
$wsImportance.Range("O12").Value = 13.9187
$wsImportance.Range("P12").Value = 0.2665119
$wsImportance.Range("Q12").Value = 4.7249813
$wsImportance.Range("R12").Value = 0.035154

$wsImportance.Range("R12").Select()

# --- "optimal models" sheet: append row 12 -----------------------------
$wsModels.Range("A12").Value = "SBFT"
$wsModels.Range("B12").Value = 2
$wsModels.Range("C12").Value = 0.75
$wsModels.Range("D12").Value = 0.006
$wsModels.Range("E12").Value = 2300
$wsModels.Range("F12").Value = 0.9978
$wsModels.Range("G12").Value = 0.9878
$wsModels.Range("H12").Value = 0.9972
$wsModels.Range("I12").Value = 0.028884

$wsModels.Activate()
$wsModels.Range("H13").Select()
